# Update "想去人数" (want-to-go count) figures and one cover image URL,
# mirroring the source refresh captured at commit 456a3b4.
# Sheet "展览" (exhibitions) and sheet "全部类型" (all types) carry the
# same rows but picked up slightly different "want to go" counts for
# row 10 when the two sheets were scraped.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Cells.Item(2, 6).Value  = 201
    $ws.Cells.Item(3, 6).Value  = 1027
    $ws.Cells.Item(5, 6).Value  = 363
    $ws.Cells.Item(6, 6).Value  = 4566
    $ws.Cells.Item(8, 6).Value  = 373
    $ws.Cells.Item(9, 6).Value  = 1327

    if ($sheetName -eq "展览") {
        $ws.Cells.Item(10, 6).Value = 878
    } else {
        $ws.Cells.Item(10, 6).Value = 879
    }
    $ws.Cells.Item(10, 9).Value = "//i2.hdslb.com/bfs/openplatform/202402/oM49o66R1708334630235.jpeg"

    $ws.Cells.Item(12, 6).Value = 935
    $ws.Cells.Item(14, 6).Value = 521
    $ws.Cells.Item(15, 6).Value = 56
    $ws.Cells.Item(16, 6).Value = 245
}
